$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.226.37'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '1.858.57'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7120'
$ws.Range("E5").Value = '  +1.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.46'
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3100'
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07752'
$ws.Range("E9").Value = '  -0.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.84'
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07808'
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("D12").Value = '1.865.37'
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '92.08'
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.094'
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6873'
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.500'
$ws.Range("E16").Value = '  +2.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008432'
$ws.Range("E17").Value = '  +2.27%  '
$ws.Range("D18").Value = '29.226.95'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '250.12'
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").Value = '2.107.46'
$ws.Range("E20").Value = '  -2.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.83'
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9997'
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.521'
$ws.Range("E23").Value = '  -1.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9999'
$ws.Range("E24").Value = '  -0.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1545'
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.47'
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.859'
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.50'
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.561'
$ws.Range("E29").Value = '  +4.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.244'
$ws.Range("E30").Value = '  -0.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.225'
$ws.Range("E31").Value = '  -0.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.198'
$ws.Range("E32").Value = '  -1.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05204'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7597'
$ws.Range("E34").Value = '  +1.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.843'
$ws.Range("E35").Value = '  -2.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.165'
$ws.Range("E36").Value = '  +0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.706'
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01860'
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("D39").Value = '1.218.74'
$ws.Range("E39").Value = '  -3.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.728'
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8962'
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.87'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9997'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.576'
$ws.Range("D45").Value = '2.002.30'
$ws.Range("E45").Value = '  -4.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5179'
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '64.61'
$ws.Range("E47").Value = '  -9.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.510'
$ws.Range("E48").Value = '  +1.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000122'
$ws.Range("E49").Value = '  -7.35%  '
$ws.Range("E50").Value = '  -2.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.006'
$ws.Range("E51").Value = '  +0.39%  '
